# Fill in the previously-empty row 41 of the work-time log with a new
# entry (Date / Start Time / End Time). The "Work Time" column (D) is a
# shared formula (=ABS(C-B)) that recalculates automatically, and the
# grand-total cell D50 (=SUM(Table1[Work Time])*24) picks up the change
# too, so we only need to touch A41:C41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serial 43069 -> 2017-11-30
$ws.Range("A41").Value = 43069
# Start time 7:40 AM
$ws.Range("B41").Value = 0.31944444444444448
# End time 8:40 AM
$ws.Range("C41").Value = 0.3611111111111111
